# DOMA-3100 add formatter convert to number for some colomns
#
# The ticket analytics export template uses text placeholders such as
# "{d.tickets[i].processing}" in row 2 (current period) and row 3 (previous
# period), columns B..G. This change appends the Carbone ":formatN()"
# formatter to those placeholders so the rendered report engine converts
# the values to numbers, and switches the cell number format of those
# columns from Text ("@") to a plain integer format ("0") so the numeric
# result displays correctly.
#
# Column A (the "address" placeholders) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: {d.tickets[i].*} placeholders -----------------------------
$ws.Range("B2").Value = "{d.tickets[i].processing:formatN()}"
$ws.Range("C2").Value = "{d.tickets[i].completed:formatN()}"
$ws.Range("D2").Value = "{d.tickets[i].canceled:formatN()}"
$ws.Range("E2").Value = "{d.tickets[i].deferred:formatN()}"
$ws.Range("F2").Value = "{d.tickets[i].closed:formatN()}"
$ws.Range("G2").Value = "{d.tickets[i].new_or_reopened:formatN()}"

# --- Row 3: {d.tickets[i+1].*} placeholders ----------------------------
$ws.Range("B3").Value = "{d.tickets[i+1].processing:formatN()}"
$ws.Range("C3").Value = "{d.tickets[i+1].completed:formatN()}"
$ws.Range("D3").Value = "{d.tickets[i+1].canceled:formatN()}"
$ws.Range("E3").Value = "{d.tickets[i+1].deferred:formatN()}"
$ws.Range("F3").Value = "{d.tickets[i+1].closed:formatN()}"
$ws.Range("G3").Value = "{d.tickets[i+1].new_or_reopened:formatN()}"

# --- Switch number counts columns (B:G) on rows 2 and 3 to a numeric format
$ws.Range("B2:G2").NumberFormat = "0"
$ws.Range("B3:G3").NumberFormat = "0"
